# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume/coin-order update described in the commit
# "Updated cryptos list on Thu Oct 31 15:13:00 UTC 2024 with GitHub Actions"
#
# Every changed cell in this sheet holds plain text (coin names, coinranking.com
# links, "70.928.68"-style prices, "  -1.93%  "-style deltas). A leading apostrophe
# is used on every write so Excel's COM layer stores the literal text instead of
# auto-coercing number-looking values (e.g. "1.00" -> 1, "147.00" -> 147) -- this
# mirrors how the source file keeps these as inline/shared strings, not numerics,
# and avoids introducing a new (Text) number-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a hashtable of column letter -> new text value
$updates = @(
    @{ Row = 2; Cols = @{ 'D'='70.928.68'; 'E'='  -1.93%  ' } }
    @{ Row = 3; Cols = @{ 'D'='2.565.44'; 'E'='  -5.32%  ' } }
    @{ Row = 4; Cols = @{ 'D'='0.999'; 'E'='  -0.15%  ' } }
    @{ Row = 5; Cols = @{ 'D'='580.06'; 'E'='  -3.43%  ' } }
    @{ Row = 6; Cols = @{ 'D'='171.33'; 'E'='  -2.31%  ' } }
    @{ Row = 7; Cols = @{ 'E'='  -0.09%  ' } }
    @{ Row = 8; Cols = @{ 'D'='0.512'; 'E'='  -2.50%  ' } }
    @{ Row = 9; Cols = @{ 'B'='LidoStakedEther'; 'C'='https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'; 'D'='2.564.41'; 'E'='  -5.32%  ' } }
    @{ Row = 10; Cols = @{ 'B'='Dogecoin'; 'C'='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; 'D'='0.167'; 'E'='  -1.24%  ' } }
    @{ Row = 11; Cols = @{ 'D'='0.169'; 'E'='  -0.22%  ' } }
    @{ Row = 12; Cols = @{ 'D'='0.352'; 'E'='  -1.07%  ' } }
    @{ Row = 13; Cols = @{ 'D'='4.86'; 'E'='  -3.06%  ' } }
    @{ Row = 14; Cols = @{ 'D'='3.054.82'; 'E'='  -4.68%  ' } }
    @{ Row = 15; Cols = @{ 'D'='0.0000184'; 'E'='  -0.93%  ' } }
    @{ Row = 16; Cols = @{ 'D'='70.677.15'; 'E'='  -2.01%  ' } }
    @{ Row = 17; Cols = @{ 'D'='25.29'; 'E'='  -4.11%  ' } }
    @{ Row = 18; Cols = @{ 'D'='2.551.03'; 'E'='  -5.80%  ' } }
    @{ Row = 19; Cols = @{ 'D'='11.77'; 'E'='  -4.54%  ' } }
    @{ Row = 20; Cols = @{ 'D'='7.69'; 'E'='  -6.39%  ' } }
    @{ Row = 21; Cols = @{ 'D'='365.26'; 'E'='  -2.46%  ' } }
    @{ Row = 22; Cols = @{ 'D'='4.02'; 'E'='  -4.27%  ' } }
    @{ Row = 23; Cols = @{ 'D'='2.01'; 'E'='  -0.65%  ' } }
    @{ Row = 24; Cols = @{ 'D'='1.00'; 'E'='  -0.01%  ' } }
    @{ Row = 25; Cols = @{ 'D'='70.22'; 'E'='  -2.99%  ' } }
    @{ Row = 26; Cols = @{ 'D'='4.17'; 'E'='  -4.90%  ' } }
    @{ Row = 27; Cols = @{ 'D'='9.32'; 'E'='  -5.32%  ' } }
    @{ Row = 28; Cols = @{ 'D'='2.682.74'; 'E'='  -5.82%  ' } }
    @{ Row = 29; Cols = @{ 'D'='0.999'; 'E'='  +0.03%  ' } }
    @{ Row = 30; Cols = @{ 'D'='0.0₃0936'; 'E'='  -5.49%  ' } }
    @{ Row = 31; Cols = @{ 'D'='7.84'; 'E'='  -3.74%  ' } }
    @{ Row = 32; Cols = @{ 'B'='Fetch.AI'; 'C'='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; 'D'='1.32'; 'E'='  -0.12%  ' } }
    @{ Row = 33; Cols = @{ 'B'='Bittensor'; 'C'='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; 'D'='484.53'; 'E'='  -4.61%  ' } }
    @{ Row = 34; Cols = @{ 'E'='  -2.97%  ' } }
    @{ Row = 35; Cols = @{ 'E'='  -0.11%  ' } }
    @{ Row = 36; Cols = @{ 'D'='157.27'; 'E'='  -4.10%  ' } }
    @{ Row = 37; Cols = @{ 'D'='0.113'; 'E'='  +4.95%  ' } }
    @{ Row = 38; Cols = @{ 'D'='18.85'; 'E'='  -4.47%  ' } }
    @{ Row = 39; Cols = @{ 'D'='18.84'; 'E'='  -1.51%  ' } }
    @{ Row = 40; Cols = @{ 'D'='1.34'; 'E'='  -4.27%  ' } }
    @{ Row = 41; Cols = @{ 'E'='  +0.00%  ' } }
    @{ Row = 42; Cols = @{ 'E'='  -5.79%  ' } }
    @{ Row = 43; Cols = @{ 'D'='2.49'; 'E'='  -3.08%  ' } }
    @{ Row = 44; Cols = @{ 'D'='4.80'; 'E'='  -5.46%  ' } }
    @{ Row = 45; Cols = @{ 'D'='0.322'; 'E'='  -3.94%  ' } }
    @{ Row = 46; Cols = @{ 'D'='38.54'; 'E'='  -2.56%  ' } }
    @{ Row = 47; Cols = @{ 'D'='147.00'; 'E'='  -6.94%  ' } }
    @{ Row = 48; Cols = @{ 'D'='3.59'; 'E'='  -5.18%  ' } }
    @{ Row = 49; Cols = @{ 'D'='0.535'; 'E'='  -5.80%  ' } }
    @{ Row = 50; Cols = @{ 'D'='1.65'; 'E'='  -7.85%  ' } }
    @{ Row = 51; Cols = @{ 'D'='0.597'; 'E'='  -1.87%  ' } }
)

foreach ($u in $updates) {
    foreach ($col in $u.Cols.Keys) {
        # Leading apostrophe forces text entry, matching the original string cells
        $ws.Range("$col$($u.Row)").Value = "'" + $u.Cols[$col]
    }
}
